$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1. Remove the "Group Reflection - Brandon" heading paragraph, merging
#    its following paragraph ("After completing Assessment 2, ...") up
#    so it becomes the new first paragraph of the body.
# ---------------------------------------------------------------------
$heading = $d.Paragraphs(1)
$headingRange = $d.Range($heading.Range.Start, $heading.Range.End)
$headingRange.Delete() | Out-Null

# ---------------------------------------------------------------------
# 2. The paragraph that now leads the document ("After completing
#    Assessment 2 ...") should carry the Normal style rather than
#    inheriting Heading 1.
# ---------------------------------------------------------------------
$d.Paragraphs(1).Style = "Normal"

# Re-assert "Normal" explicitly on the other two body paragraphs too
# (they already render as Normal, this just makes the style explicit).
$d.Paragraphs(2).Style = "Normal"
$d.Paragraphs(3).Style = "Normal"

# ---------------------------------------------------------------------
# 3. Section / page-setup touch-ups.
# ---------------------------------------------------------------------
$section = $d.Sections(1)
$pageSetup = $section.PageSetup
$pageSetup.HeaderDistance = 0
$pageSetup.FooterDistance = 0
$section.ProtectedForForms = $false
